{"js": "// Apply the cover-letter revisions described by the diff.\n// Strategy: locate each affected stretch of text with Body.search()\n// (exact, case-sensitive match spanning the existing runs) and replace it\n// in place with Range.insertText(..., \"Replace\"). This mirrors how a human\n// editor would use Word's Find & Replace and keeps paragraph / run\n// structure outside the edited spans untouched.\n\nconst body = context.document.body;\n\nasync function replaceOnce(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"Coupland, and I am trying...\" -> \"...I\u2019m trying...\"\nawait replaceOnce(\n  \"Coupland, and I am trying to persuade you to hire me.\",\n  \"Coupland, and I\\u2019m trying to persuade you to hire me.\"\n);\n\n// 2) \"...but I am very proficient in...\" -> \"...but I\u2019m very proficient with...\"\n//    and drop the trailing \"I also try to make time for volunteer work.\" sentence\n//    (that idea moves into the next paragraph).\nawait replaceOnce(\n  \"but I am very proficient in the entire alphabet soup of tools and languages. I have delivered software of all types including applications, embedded systems, APIs, and services, and some simple Android apps. I also try to make time for volunteer work.\",\n  \"but I\\u2019m very proficient with the entire alphabet soup of tools and languages. I have delivered software of all types including applications, embedded systems, APIs, and services, as well as some simple Android apps. \"\n);\n\n// 3) Rework the \"step away from the keyboard\" paragraph: new lead-in about\n//    volunteering vs. hobbies, \"robotics\" -> \"robotic\", insert \"some\" before\n//    RaspberryPi, and add the \"noir fiction\" aside.\nawait replaceOnce(\n  \"When I need to step away from the keyboard and recharge, I find something to build, go camping, or go to concerts with my sister. I like animals, tinkering in my garage, and discovering new (or new to me) bands. I also write code for fun. Lately I have been working on \",\n  \"When I need to step away from the keyboard and recharge, I try to find something useful to keep me busy. Sometimes that means volunteering. Other times it\\u2019s something closer to home like finding something to build, going camping, or going to concerts with my sister. I like animals, tinkering in my garage, and discovering new (or new to me) bands. I write noir fiction. I also write code. For fun. Lately I have been working on some \"\n);\nawait replaceOnce(\", robotics,\", \", robotic,\");\n\n// 4) Alpaca description + history rewrite.\nawait replaceOnce(\n  \"is Alpaca, a tool to manage and manipulate images, videos, and audio. Alpaca started as \",\n  \"is Alpaca, a multi-media management tool. Alpaca started as \"\n);\nawait replaceOnce(\n  \" application. Software for a custom interface to\",\n  \" application - a custom interface to\"\n);\nawait replaceOnce(\n  \", but I changed directions and now it is a desktop application. Version 1.0 will ship once there are no new tools or techniques for me to learn. \",\n  \", but it\\u2019s morphed into a desktop application over the years. Version 1.0 will ship once there are no new tools or techniques for me to learn. \"\n);\n\n// 5) \"have not filled this position yet\" -> \"haven\u2019t filled this position\"\nawait replaceOnce(\n  \"I sincerely hope you have not filled this position yet and that you have time to look at my resume.\",\n  \"I sincerely hope you haven\\u2019t filled this position and that you have time to look at my resume.\"\n);\n\n// 6) \"very aware\" -> \"well aware\"\nawait replaceOnce(\n  \"I lived in Alaska for 3 years, so I am very aware of the time zone difference.\",\n  \"I lived in Alaska for 3 years, so I am well aware of the time zone difference.\"\n);\n", "ps1": "# Apply the cover-letter revisions described by the diff.\n# Strategy: use Find/Replace (Range.Find.Execute) against $d.Content for each\n# affected stretch of text. Each \"find\" string is long/specific enough to be\n# unique in the document, so this mirrors how a human editor would use\n# Word's Find & Replace dialog (Replace All) without disturbing unrelated\n# paragraphs or runs.\n\n$d = $word.ActiveDocument\n\nfunction ReplaceOnce([string]$findText, [string]$replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #          MatchAllWordForms, Forward, Wrap(1=wdFindContinue), Format,\n    #          ReplaceWith, Replace(2=wdReplaceAll))\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1) \"Coupland, and I am trying...\" -> \"...I'm trying...\"\nReplaceOnce \"Coupland, and I am trying to persuade you to hire me.\" \"Coupland, and I\u2019m trying to persuade you to hire me.\"\n\n# 2) \"...but I am very proficient in...\" -> \"...but I'm very proficient with...\"\n#    and drop the trailing \"I also try to make time for volunteer work.\" sentence\n#    (that idea moves into the next paragraph).\nReplaceOnce \"but I am very proficient in the entire alphabet soup of tools and languages. I have delivered software of all types including applications, embedded systems, APIs, and services, and some simple Android apps. I also try to make time for volunteer work.\" \"but I\u2019m very proficient with the entire alphabet soup of tools and languages. I have delivered software of all types including applications, embedded systems, APIs, and services, as well as some simple Android apps. \"\n\n# 3) Rework the \"step away from the keyboard\" paragraph: new lead-in about\n#    volunteering vs. hobbies, \"robotics\" -> \"robotic\", insert \"some\" before\n#    RaspberryPi, and add the \"noir fiction\" aside.\nReplaceOnce \"When I need to step away from the keyboard and recharge, I find something to build, go camping, or go to concerts with my sister. I like animals, tinkering in my garage, and discovering new (or new to me) bands. I also write code for fun. Lately I have been working on \" \"When I need to step away from the keyboard and recharge, I try to find something useful to keep me busy. Sometimes that means volunteering. Other times it\u2019s something closer to home like finding something to build, going camping, or going to concerts with my sister. I like animals, tinkering in my garage, and discovering new (or new to me) bands. I write noir fiction. I also write code. For fun. Lately I have been working on some \"\nReplaceOnce \", robotics,\" \", robotic,\"\n\n# 4) Alpaca description + history rewrite.\nReplaceOnce \"is Alpaca, a tool to manage and manipulate images, videos, and audio. Alpaca started as \" \"is Alpaca, a multi-media management tool. Alpaca started as \"\nReplaceOnce \" application. Software for a custom interface to\" \" application - a custom interface to\"\nReplaceOnce \", but I changed directions and now it is a desktop application. Version 1.0 will ship once there are no new tools or techniques for me to learn. \" \", but it\u2019s morphed into a desktop application over the years. Version 1.0 will ship once there are no new tools or techniques for me to learn. \"\n\n# 5) \"have not filled this position yet\" -> \"haven't filled this position\"\nReplaceOnce \"I sincerely hope you have not filled this position yet and that you have time to look at my resume.\" \"I sincerely hope you haven\u2019t filled this position and that you have time to look at my resume.\"\n\n# 6) \"very aware\" -> \"well aware\"\nReplaceOnce \"I lived in Alaska for 3 years, so I am very aware of the time zone difference.\" \"I lived in Alaska for 3 years, so I am well aware of the time zone difference.\"\n"}
